$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("forest")

$years = 1991..2020
$values = @(
    [double]"9.2437438999999996E-2",
    [double]"9.2759037000000003E-2",
    [double]"9.3386048999999999E-2",
    [double]"9.3757749000000001E-2",
    [double]"9.4563718000000005E-2",
    [double]"9.5301708999999998E-2",
    [double]"9.5933589999999999E-2",
    [double]"9.6374731000000005E-2",
    [double]"9.6559106000000006E-2",
    [double]"9.6564838E-2",
    [double]"9.6685108000000006E-2",
    [double]"9.6855833000000002E-2",
    [double]"9.7477267000000006E-2",
    [double]"9.8103177999999999E-2",
    [double]"9.9019310999999999E-2",
    0.100250851,
    0.101453266,
    [double]"0.10213665199999999",
    0.102838016,
    [double]"0.10339466999999999",
    0.104658512,
    0.105491901,
    0.106636957,
    0.107014838,
    0.107658669,
    0.108301522,
    [double]"0.10982653100000001",
    [double]"0.11090301700000001",
    [double]"0.11170989000000001",
    0.112605549
)

$startRow = 122
for ($i = 0; $i -lt $years.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $years[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Range("A122:B151").VerticalAlignment = -4108

$ws.Range("A122:A151").Select()

$excel.ActiveWindow.ScrollRow = 134
